# Fix the logic error in the loan-installment formulas ("solved error in
# logic for the functions"):
#  - Row 5 ("Cuota durante los estudios") formulas previously used
#    inconsistent hard-coded percentages (0.3/0.6/none) and exponent -1.5;
#    corrected to consistently multiply by the term (C3, D3, ...) with a
#    fixed 0.3 factor and use exponent -1.
#  - Row 6 ("Cuota después de los estudios") formulas previously used
#    inconsistent hard-coded percentages/rates (0.7 vs 0.4/0.0099); corrected
#    to consistently multiply by the term and use a fixed 0.3 / 0.0115 rate
#    with exponent -1.5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: "Cuota durante los estudios"
$ws.Range("C5").Formula = "= (0.3 * C2 * C3 * 0.0115) / (1 - (1 + 0.0115)^(-1 * (C3/2) * 12))"
$ws.Range("D5").Formula = "= (0.3 * D2 * D3 * 0.0115) / (1 - (1 + 0.0115)^(-1 * (D3/2) * 12))"
$ws.Range("E5").Formula = "= (0.3 * E2 * E3 * 0.0115) / (1 - (1 + 0.0115)^(-1 * (E3/2) * 12))"
$ws.Range("F5").Formula = "= (0.3 * F2 * F3 * 0.0115) / (1 - (1 + 0.0115)^(-1 * (F3/2) * 12))"
$ws.Range("G5").Formula = "= (0.3 * G2 * G3 * 0.0115) / (1 - (1 + 0.0115)^(-1 * (G3/2) * 12))"
$ws.Range("H5").Formula = "= (0.3 * H2 * H3 * 0.0115) / (1 - (1 + 0.0115)^(-1 * (H3/2) * 12))"

# Row 6: "Cuota después de los estudios"
$ws.Range("C6").Formula = "= (0.3 * C2 * C3 * 0.0115) / (1 - (1 + 0.0115)^(-1.5 * (C3/2) * 12))"
$ws.Range("D6").Formula = "= (0.3 * D2 * D3 * 0.0115) / (1 - (1 + 0.0115)^(-1.5 * (D3/2) * 12))"
$ws.Range("E6").Formula = "= (0.3 * E2 * E3 * 0.0115) / (1 - (1 + 0.0115)^(-1.5 * (E3/2) * 12))"
$ws.Range("F6").Formula = "= (0.3 * F2 * F3 * 0.0115) / (1 - (1 + 0.0115)^(-1.5 * (F3/2) * 12))"

$wb.Save()
